$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.058.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -4.28%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.240.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -7.59%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'594.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.58%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'153.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -11.00%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.11%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'3.231.41"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -7.69%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.548"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -10.03%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  -9.98%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'6.69"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -7.94%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.503"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -14.14%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'39.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -14.48%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.0000249"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -9.72%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'3.760.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -7.65%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'67.139.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.19%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.235.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -7.59%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  -4.58%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'7.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -13.47%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'533.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -13.19%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'15.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -13.44%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  -12.79%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'7.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -13.60%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'13.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -10.46%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'86.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -12.76%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  -0.03%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  -14.08%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  -13.53%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'8.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -8.96%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'29.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -12.57%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'2.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -11.64%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'1.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -9.59%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").Value = "'Filecoin"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'6.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -17.91%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").Value = "'Bittensor"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'537.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -14.94%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'5.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -14.97%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +0.14%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'53.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -6.12%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  -12.17%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").Value = "'VeChain"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'0.0427"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -11.05%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = "'Cosmos"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'9.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -12.59%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  -11.81%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'2.948.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -12.12%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  -23.13%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'  -13.40%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.0" + [char]0x2083 + "0596"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -18.82%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D47").Value = "'26.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -16.50%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = "'Fetch.AI"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'2.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -16.01%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "'USDe"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.15%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  -11.71%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'122.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -7.83%  "
$ws.Range("E51").Style = "Normal"
